$d = $word.ActiveDocument

$d.Content.Find.Execute("43+28=71", $true, $false, $false, $false, $false, $true, 1, $false, "81-47=34", 2) | Out-Null
$d.Content.Find.Execute("30+0=30", $true, $false, $false, $false, $false, $true, 1, $false, "72-7=65", 2) | Out-Null
$d.Content.Find.Execute("28+69=97", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 2) | Out-Null
$d.Content.Find.Execute("82-9=73", $true, $false, $false, $false, $false, $true, 1, $false, "13+76=89", 2) | Out-Null
$d.Content.Find.Execute("5+20=25", $true, $false, $false, $false, $false, $true, 1, $false, "19+5=24", 2) | Out-Null
$d.Content.Find.Execute("46+18=64", $true, $false, $false, $false, $false, $true, 1, $false, "99-94=5", 2) | Out-Null
$d.Content.Find.Execute("24-6=18", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=2", 2) | Out-Null
$d.Content.Find.Execute("84-7=77", $true, $false, $false, $false, $false, $true, 1, $false, "54+45=99", 2) | Out-Null
$d.Content.Find.Execute("55-3=52", $true, $false, $false, $false, $false, $true, 1, $false, "74-23=51", 2) | Out-Null
$d.Content.Find.Execute("92-36=56", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=41", 2) | Out-Null
$d.Content.Find.Execute("3+65=68", $true, $false, $false, $false, $false, $true, 1, $false, "27-23=4", 2) | Out-Null
$d.Content.Find.Execute("99-26=73", $true, $false, $false, $false, $false, $true, 1, $false, "69+2=71", 2) | Out-Null
$d.Content.Find.Execute("26+24=50", $true, $false, $false, $false, $false, $true, 1, $false, "9+25=34", 2) | Out-Null
$d.Content.Find.Execute("83-51=32", $true, $false, $false, $false, $false, $true, 1, $false, "15+64=79", 2) | Out-Null
$d.Content.Find.Execute("99-91=8", $true, $false, $false, $false, $false, $true, 1, $false, "5-4=1", 2) | Out-Null
$d.Content.Find.Execute("83-6=77", $true, $false, $false, $false, $false, $true, 1, $false, "29-12=17", 2) | Out-Null
$d.Content.Find.Execute("85-55=30", $true, $false, $false, $false, $false, $true, 1, $false, "74-37=37", 2) | Out-Null
$d.Content.Find.Execute("54-46=8", $true, $false, $false, $false, $false, $true, 1, $false, "88-66=22", 2) | Out-Null
$d.Content.Find.Execute("94-70=24", $true, $false, $false, $false, $false, $true, 1, $false, "2+26=28", 2) | Out-Null
$d.Content.Find.Execute("13-0=13", $true, $false, $false, $false, $false, $true, 1, $false, "81-36=45", 2) | Out-Null
$d.Content.Find.Execute("3+18=21", $true, $false, $false, $false, $false, $true, 1, $false, "50-18=32", 2) | Out-Null
$d.Content.Find.Execute("64+28=92", $true, $false, $false, $false, $false, $true, 1, $false, "37-7=30", 2) | Out-Null
$d.Content.Find.Execute("70+27=97", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=84", 2) | Out-Null
$d.Content.Find.Execute("38+23=61", $true, $false, $false, $false, $false, $true, 1, $false, "85+9=94", 2) | Out-Null
$d.Content.Find.Execute("67-2=65", $true, $false, $false, $false, $false, $true, 1, $false, "73+3=76", 2) | Out-Null
$d.Content.Find.Execute("72-49=23", $true, $false, $false, $false, $false, $true, 1, $false, "62+34=96", 2) | Out-Null
$d.Content.Find.Execute("37+47=84", $true, $false, $false, $false, $false, $true, 1, $false, "68-3=65", 2) | Out-Null
$d.Content.Find.Execute("63-50=13", $true, $false, $false, $false, $false, $true, 1, $false, "24+24=48", 2) | Out-Null
$d.Content.Find.Execute("11+68=79", $true, $false, $false, $false, $false, $true, 1, $false, "45+38=83", 2) | Out-Null
$d.Content.Find.Execute("30-12=18", $true, $false, $false, $false, $false, $true, 1, $false, "41-29=12", 2) | Out-Null
$d.Content.Find.Execute("53-30=23", $true, $false, $false, $false, $false, $true, 1, $false, "33-8=25", 2) | Out-Null
$d.Content.Find.Execute("56+16=72", $true, $false, $false, $false, $false, $true, 1, $false, "39+4=43", 2) | Out-Null
$d.Content.Find.Execute("63-38=25", $true, $false, $false, $false, $false, $true, 1, $false, "66-18=48", 2) | Out-Null
$d.Content.Find.Execute("42+38=80", $true, $false, $false, $false, $false, $true, 1, $false, "27+58=85", 2) | Out-Null
$d.Content.Find.Execute("12+82=94", $true, $false, $false, $false, $false, $true, 1, $false, "68-12=56", 2) | Out-Null
$d.Content.Find.Execute("94-36=58", $true, $false, $false, $false, $false, $true, 1, $false, "23+41=64", 2) | Out-Null
$d.Content.Find.Execute("95-65=30", $true, $false, $false, $false, $false, $true, 1, $false, "34+30=64", 2) | Out-Null
$d.Content.Find.Execute("14+77=91", $true, $false, $false, $false, $false, $true, 1, $false, "63-16=47", 2) | Out-Null
$d.Content.Find.Execute("80-36=44", $true, $false, $false, $false, $false, $true, 1, $false, "56+18=74", 2) | Out-Null
$d.Content.Find.Execute("4+83=87", $true, $false, $false, $false, $false, $true, 1, $false, "70+12=82", 2) | Out-Null
$d.Content.Find.Execute("11+75=86", $true, $false, $false, $false, $false, $true, 1, $false, "34-27=7", 2) | Out-Null
$d.Content.Find.Execute("30+53=83", $true, $false, $false, $false, $false, $true, 1, $false, "43-4=39", 2) | Out-Null
$d.Content.Find.Execute("69-30=39", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=39", 2) | Out-Null
$d.Content.Find.Execute("36+20=56", $true, $false, $false, $false, $false, $true, 1, $false, "27+28=55", 2) | Out-Null
$d.Content.Find.Execute("44-35=9", $true, $false, $false, $false, $false, $true, 1, $false, "49+33=82", 2) | Out-Null
$d.Content.Find.Execute("5+86=91", $true, $false, $false, $false, $false, $true, 1, $false, "56-3=53", 2) | Out-Null
$d.Content.Find.Execute("16+11=27", $true, $false, $false, $false, $false, $true, 1, $false, "41+26=67", 2) | Out-Null
$d.Content.Find.Execute("43+45=88", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=60", 2) | Out-Null
$d.Content.Find.Execute("35-18=17", $true, $false, $false, $false, $false, $true, 1, $false, "39-34=5", 2) | Out-Null
$d.Content.Find.Execute("86-77=9", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=24", 2) | Out-Null
$d.Content.Find.Execute("23+31=54", $true, $false, $false, $false, $false, $true, 1, $false, "51-50=1", 2) | Out-Null
$d.Content.Find.Execute("42-9=33", $true, $false, $false, $false, $false, $true, 1, $false, "44+25=69", 2) | Out-Null
$d.Content.Find.Execute("7+30=37", $true, $false, $false, $false, $false, $true, 1, $false, "96-15=81", 2) | Out-Null
$d.Content.Find.Execute("87-64=23", $true, $false, $false, $false, $false, $true, 1, $false, "29+36=65", 2) | Out-Null
$d.Content.Find.Execute("3+60=63", $true, $false, $false, $false, $false, $true, 1, $false, "43+35=78", 2) | Out-Null
$d.Content.Find.Execute("56+34=90", $true, $false, $false, $false, $false, $true, 1, $false, "44-2=42", 2) | Out-Null
$d.Content.Find.Execute("84-21=63", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=38", 2) | Out-Null
$d.Content.Find.Execute("26+15=41", $true, $false, $false, $false, $false, $true, 1, $false, "74-70=4", 2) | Out-Null
$d.Content.Find.Execute("51-25=26", $true, $false, $false, $false, $false, $true, 1, $false, "54-22=32", 2) | Out-Null
$d.Content.Find.Execute("58+28=86", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=18", 2) | Out-Null
$d.Content.Find.Execute("14+7=21", $true, $false, $false, $false, $false, $true, 1, $false, "96-10=86", 2) | Out-Null
$d.Content.Find.Execute("11+4=15", $true, $false, $false, $false, $false, $true, 1, $false, "67-16=51", 2) | Out-Null
$d.Content.Find.Execute("20+32=52", $true, $false, $false, $false, $false, $true, 1, $false, "45+35=80", 2) | Out-Null
$d.Content.Find.Execute("56-51=5", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=85", 2) | Out-Null
$d.Content.Find.Execute("60-39=21", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=52", 2) | Out-Null
$d.Content.Find.Execute("69-54=15", $true, $false, $false, $false, $false, $true, 1, $false, "56-15=41", 2) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "2+62=64", 2) | Out-Null
$d.Content.Find.Execute("11+70=81", $true, $false, $false, $false, $false, $true, 1, $false, "49-30=19", 2) | Out-Null
$d.Content.Find.Execute("81-26=55", $true, $false, $false, $false, $false, $true, 1, $false, "99-86=13", 2) | Out-Null
$d.Content.Find.Execute("80-5=75", $true, $false, $false, $false, $false, $true, 1, $false, "1+28=29", 2) | Out-Null
$d.Content.Find.Execute("40+59=99", $true, $false, $false, $false, $false, $true, 1, $false, "22+65=87", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $true, $false, $false, $false, $false, $true, 1, $false, "8+1=9", 2) | Out-Null
$d.Content.Find.Execute("52-48=4", $true, $false, $false, $false, $false, $true, 1, $false, "52-26=26", 2) | Out-Null
$d.Content.Find.Execute("29-4=25", $true, $false, $false, $false, $false, $true, 1, $false, "15+60=75", 2) | Out-Null
$d.Content.Find.Execute("91-44=47", $true, $false, $false, $false, $false, $true, 1, $false, "52+31=83", 2) | Out-Null
$d.Content.Find.Execute("19+4=23", $true, $false, $false, $false, $false, $true, 1, $false, "15+59=74", 2) | Out-Null
$d.Content.Find.Execute("0+8=8", $true, $false, $false, $false, $false, $true, 1, $false, "56-39=17", 2) | Out-Null
$d.Content.Find.Execute("46+24=70", $true, $false, $false, $false, $false, $true, 1, $false, "0+3=3", 2) | Out-Null
$d.Content.Find.Execute("79+11=90", $true, $false, $false, $false, $false, $true, 1, $false, "83+3=86", 2) | Out-Null
$d.Content.Find.Execute("21+37=58", $true, $false, $false, $false, $false, $true, 1, $false, "66+12=78", 2) | Out-Null
$d.Content.Find.Execute("95-52=43", $true, $false, $false, $false, $false, $true, 1, $false, "52+3=55", 2) | Out-Null
$d.Content.Find.Execute("43+34=77", $true, $false, $false, $false, $false, $true, 1, $false, "66-13=53", 2) | Out-Null
$d.Content.Find.Execute("67-37=30", $true, $false, $false, $false, $false, $true, 1, $false, "22+0=22", 2) | Out-Null
$d.Content.Find.Execute("85-81=4", $true, $false, $false, $false, $false, $true, 1, $false, "17+47=64", 2) | Out-Null
$d.Content.Find.Execute("41+21=62", $true, $false, $false, $false, $false, $true, 1, $false, "54-53=1", 2) | Out-Null
$d.Content.Find.Execute("3+27=30", $true, $false, $false, $false, $false, $true, 1, $false, "97-39=58", 2) | Out-Null
$d.Content.Find.Execute("23+34=57", $true, $false, $false, $false, $false, $true, 1, $false, "93-5=88", 2) | Out-Null
$d.Content.Find.Execute("76-47=29", $true, $false, $false, $false, $false, $true, 1, $false, "94-8=86", 2) | Out-Null
$d.Content.Find.Execute("30-13=17", $true, $false, $false, $false, $false, $true, 1, $false, "85-68=17", 2) | Out-Null
$d.Content.Find.Execute("85-27=58", $true, $false, $false, $false, $false, $true, 1, $false, "86-7=79", 2) | Out-Null
$d.Content.Find.Execute("19+79=98", $true, $false, $false, $false, $false, $true, 1, $false, "16+5=21", 2) | Out-Null
$d.Content.Find.Execute("6+10=16", $true, $false, $false, $false, $false, $true, 1, $false, "38+51=89", 2) | Out-Null
$d.Content.Find.Execute("66+15=81", $true, $false, $false, $false, $false, $true, 1, $false, "89-1=88", 2) | Out-Null
$d.Content.Find.Execute("19+32=51", $true, $false, $false, $false, $false, $true, 1, $false, "64-52=12", 2) | Out-Null
$d.Content.Find.Execute("74+7=81", $true, $false, $false, $false, $false, $true, 1, $false, "24+74=98", 2) | Out-Null
$d.Content.Find.Execute("68-38=30", $true, $false, $false, $false, $false, $true, 1, $false, "3+42=45", 2) | Out-Null
$d.Content.Find.Execute("37+33=70", $true, $false, $false, $false, $false, $true, 1, $false, "99-42=57", 2) | Out-Null
$d.Content.Find.Execute("47-31=16", $true, $false, $false, $false, $false, $true, 1, $false, "67-0=67", 2) | Out-Null
$d.Content.Find.Execute("13+85=98", $true, $false, $false, $false, $false, $true, 1, $false, "87+3=90", 2) | Out-Null
$d.Content.Find.Execute("81-5=76", $true, $false, $false, $false, $false, $true, 1, $false, "83-10=73", 2) | Out-Null
